# Updated cryptos list on Sat Nov 16 21:37:27 UTC 2024 with GitHub Actions
# Applies the latest scraped crypto price/volume figures to the sheet,
# including three rows whose coins were re-ranked (Stellar/Cronos/
# InternetComputer(DFINITY) rotated through rows 30-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, new text, and whether the text looks like a
# plain number (e.g. "216.45"). Price/volume figures are stored as TEXT
# in this sheet (note some, like "90.978.14", aren't valid numbers at
# all because of the thousands separators) so any cell whose new value
# would otherwise be auto-recognised as a number is forced back to text
# via a leading apostrophe (classic Excel "treat as text" marker), then
# the resulting quote-prefix style is cleared so the cell's style index
# is unchanged.
$updates = @(
    @{ Cell = 'D2'; New = '90.978.14'; Numeric = $false },
    @{ Cell = 'E2'; New = '  -0.40%  '; Numeric = $false },
    @{ Cell = 'D3'; New = '3.155.19'; Numeric = $false },
    @{ Cell = 'E3'; New = '  +2.13%  '; Numeric = $false },
    @{ Cell = 'E4'; New = '  +0.24%  '; Numeric = $false },
    @{ Cell = 'D5'; New = '216.45'; Numeric = $true },
    @{ Cell = 'E5'; New = '  -0.76%  '; Numeric = $false },
    @{ Cell = 'D6'; New = '626.14'; Numeric = $true },
    @{ Cell = 'D7'; New = '1.14'; Numeric = $true },
    @{ Cell = 'E7'; New = '  +29.74%  '; Numeric = $false },
    @{ Cell = 'D8'; New = '0.369'; Numeric = $true },
    @{ Cell = 'E9'; New = '  +0.01%  '; Numeric = $false },
    @{ Cell = 'D10'; New = '3.153.86'; Numeric = $false },
    @{ Cell = 'E10'; New = '  +2.20%  '; Numeric = $false },
    @{ Cell = 'D11'; New = '0.761'; Numeric = $true },
    @{ Cell = 'E11'; New = '  +14.12%  '; Numeric = $false },
    @{ Cell = 'E12'; New = '  +6.79%  '; Numeric = $false },
    @{ Cell = 'E13'; New = '  +6.64%  '; Numeric = $false },
    @{ Cell = 'E14'; New = '  -1.12%  '; Numeric = $false },
    @{ Cell = 'D15'; New = '35.13'; Numeric = $true },
    @{ Cell = 'E15'; New = '  +6.82%  '; Numeric = $false },
    @{ Cell = 'D16'; New = '90.645.90'; Numeric = $false },
    @{ Cell = 'E16'; New = '  -0.55%  '; Numeric = $false },
    @{ Cell = 'D17'; New = '3.738.49'; Numeric = $false },
    @{ Cell = 'E17'; New = '  +2.18%  '; Numeric = $false },
    @{ Cell = 'D18'; New = '3.197.97'; Numeric = $false },
    @{ Cell = 'E19'; New = '  +8.40%  '; Numeric = $false },
    @{ Cell = 'D20'; New = '14.59'; Numeric = $true },
    @{ Cell = 'E20'; New = '  +5.89%  '; Numeric = $false },
    @{ Cell = 'D21'; New = '475.39'; Numeric = $true },
    @{ Cell = 'E21'; New = '  +9.41%  '; Numeric = $false },
    @{ Cell = 'D22'; New = '0.0000211'; Numeric = $true },
    @{ Cell = 'E22'; New = '  -3.20%  '; Numeric = $false },
    @{ Cell = 'D23'; New = '9.17'; Numeric = $true },
    @{ Cell = 'E23'; New = '  +8.27%  '; Numeric = $false },
    @{ Cell = 'D24'; New = '5.20'; Numeric = $true },
    @{ Cell = 'E24'; New = '  +1.55%  '; Numeric = $false },
    @{ Cell = 'D25'; New = '5.94'; Numeric = $true },
    @{ Cell = 'E25'; New = '  +5.98%  '; Numeric = $false },
    @{ Cell = 'D26'; New = '94.91'; Numeric = $true },
    @{ Cell = 'E26'; New = '  +13.16%  '; Numeric = $false },
    @{ Cell = 'D27'; New = '12.35'; Numeric = $true },
    @{ Cell = 'E27'; New = '  +4.79%  '; Numeric = $false },
    @{ Cell = 'D28'; New = '3.323.62'; Numeric = $false },
    @{ Cell = 'E28'; New = '  +2.03%  '; Numeric = $false },
    @{ Cell = 'E29'; New = '  -0.24%  '; Numeric = $false },
    @{ Cell = 'B30'; New = 'Stellar'; Numeric = $false },
    @{ Cell = 'C30'; New = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; Numeric = $false },
    @{ Cell = 'D30'; New = '0.219'; Numeric = $true },
    @{ Cell = 'E30'; New = '  +57.11%  '; Numeric = $false },
    @{ Cell = 'B31'; New = 'Cronos'; Numeric = $false },
    @{ Cell = 'C31'; New = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Numeric = $false },
    @{ Cell = 'D31'; New = '0.163'; Numeric = $true },
    @{ Cell = 'E31'; New = '  -1.59%  '; Numeric = $false },
    @{ Cell = 'B32'; New = 'InternetComputer(DFINITY)'; Numeric = $false },
    @{ Cell = 'C32'; New = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; Numeric = $false },
    @{ Cell = 'D32'; New = '9.33'; Numeric = $true },
    @{ Cell = 'E32'; New = '  +8.16%  '; Numeric = $false },
    @{ Cell = 'E33'; New = '  +0.06%  '; Numeric = $false },
    @{ Cell = 'D34'; New = '27.32'; Numeric = $true },
    @{ Cell = 'E34'; New = '  +18.70%  '; Numeric = $false },
    @{ Cell = 'D35'; New = '521.02'; Numeric = $true },
    @{ Cell = 'E35'; New = '  +1.25%  '; Numeric = $false },
    @{ Cell = 'E36'; New = '  +5.39%  '; Numeric = $false },
    @{ Cell = 'D37'; New = '1.94'; Numeric = $true },
    @{ Cell = 'E37'; New = '  +5.83%  '; Numeric = $false },
    @{ Cell = 'D38'; New = '7.02'; Numeric = $true },
    @{ Cell = 'E38'; New = '  +1.14%  '; Numeric = $false },
    @{ Cell = 'D39'; New = '3.62'; Numeric = $true },
    @{ Cell = 'E39'; New = '  -5.48%  '; Numeric = $false },
    @{ Cell = 'E40'; New = '  +4.12%  '; Numeric = $false },
    @{ Cell = 'D41'; New = '0.0903'; Numeric = $true },
    @{ Cell = 'E41'; New = '  +25.47%  '; Numeric = $false },
    @{ Cell = 'E42'; New = '  +16.50%  '; Numeric = $false },
    @{ Cell = 'D43'; New = '22.22'; Numeric = $true },
    @{ Cell = 'E43'; New = '  -0.41%  '; Numeric = $false },
    @{ Cell = 'E44'; New = '  -0.07%  '; Numeric = $false },
    @{ Cell = 'E45'; New = '  +6.42%  '; Numeric = $false },
    @{ Cell = 'D46'; New = '0.735'; Numeric = $true },
    @{ Cell = 'E46'; New = '  +21.28%  '; Numeric = $false },
    @{ Cell = 'E48'; New = '  +12.92%  '; Numeric = $false },
    @{ Cell = 'D49'; New = '150.40'; Numeric = $true },
    @{ Cell = 'E49'; New = '  +6.11%  '; Numeric = $false },
    @{ Cell = 'D50'; New = '1.37'; Numeric = $true },
    @{ Cell = 'E50'; New = '  +10.55%  '; Numeric = $false },
    @{ Cell = 'D51'; New = '45.37'; Numeric = $true },
    @{ Cell = 'E51'; New = '  +3.53%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $cell.Value = "'" + $u.New
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.New
    }
}
